$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G11").Value = "Fallo"
$ws.Range("H11").Value = -1
